$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 1233.317016601562
$ws.Range("C2").Value = 0.9252
$ws.Range("D2").Value = 0.913100004196167
$ws.Range("E2").Value = 1.337200045585632
$ws.Range("F2").Value = 0.6298999786376953
$ws.Range("H2").Value = 0.7916

# Row 3
$ws.Range("B3").Value = 1147.790771484375
$ws.Range("C3").Value = 0.9182
$ws.Range("D3").Value = 0.9183
$ws.Range("E3").Value = 1.043799996376038
$ws.Range("F3").Value = 0.794700026512146
$ws.Range("H3").Value = 0.8375

# Row 4
$ws.Range("B4").Value = 789.2227783203125
$ws.Range("C4").Value = 0.922
$ws.Range("D4").Value = 0.9214
$ws.Range("E4").Value = 1.018800020217896
$ws.Range("F4").Value = 0.7828999757766724
$ws.Range("H4").Value = 0.8651

# Row 5
$ws.Range("B5").Value = 841.4779052734375
$ws.Range("C5").Value = 0.8875999999999999
$ws.Range("D5").Value = 0.8985
$ws.Range("E5").Value = 0.9667999744415283
$ws.Range("F5").Value = 0.5730000138282776
$ws.Range("H5").Value = 0.6623

# Row 6
$ws.Range("B6").Value = 1145.45166015625
$ws.Range("C6").Value = 0.9062
$ws.Range("D6").Value = 0.9137
$ws.Range("E6").Value = 0.9592000246047974
$ws.Range("F6").Value = 0.6822999715805054
$ws.Range("H6").Value = 0.7963

# Row 7
$ws.Range("B7").Value = 899.5681762695312
$ws.Range("C7").Value = 0.9059
$ws.Range("D7").Value = 0.9139000177383423
$ws.Range("E7").Value = 0.9584000110626221
$ws.Range("F7").Value = 0.746399998664856
$ws.Range("H7").Value = 0.7984

# Row 8
$ws.Range("B8").Value = 996.4227294921875
$ws.Range("C8").Value = 0.8929
$ws.Range("D8").Value = 0.8972
$ws.Range("E8").Value = 0.9605000019073486
$ws.Range("F8").Value = 0.7840999960899353
$ws.Range("H8").Value = 0.6501

# Row 9
$ws.Range("B9").Value = 7053.25146484375
$ws.Range("C9").Value = 0.9089
$ws.Range("D9").Value = 0.9121
$ws.Range("E9").Value = 1.337200045585632
$ws.Range("F9").Value = 0.5730000138282776
$ws.Range("H9").Value = 5.401300000000001
